$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a value that must remain plain text even though it looks like a
# number (e.g. "557.74"), without leaving a NumberFormat override behind on
# the destination cell. The text is staged in a scratch cell (forced to
# stay text via a Text number format), then copy / paste-special carries
# only the value (and its "is text" flag, not the scratch formatting) into
# the destination; the scratch column is removed afterwards.
function Set-SafeTextValue($CellAddress, $TextValue) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $TextValue
    $scratch.Copy()
    $ws.Range($CellAddress).PasteSpecial(-4163)
    $scratch.EntireColumn.Delete()
}

$ws.Range("D2").Value = "65.699.57"
$ws.Range("E2").Value = "  -5.66%  "

$ws.Range("D3").Value = "3.268.84"
$ws.Range("E3").Value = "  -6.32%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-SafeTextValue "D5" "557.74"
$ws.Range("E5").Value = "  -3.76%  "

Set-SafeTextValue "D6" "184.57"
$ws.Range("E6").Value = "  -4.41%  "

Set-SafeTextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.06%  "

Set-SafeTextValue "D8" "0.588"
$ws.Range("E8").Value = "  -4.03%  "

$ws.Range("D9").Value = "3.263.64"
$ws.Range("E9").Value = "  -6.17%  "

Set-SafeTextValue "D10" "0.184"
$ws.Range("E10").Value = "  -9.84%  "

Set-SafeTextValue "D11" "0.582"
$ws.Range("E11").Value = "  -5.87%  "

Set-SafeTextValue "D12" "47.14"
$ws.Range("E12").Value = "  -8.53%  "

$ws.Range("E13").Value = "  -7.24%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-SafeTextValue "D14" "8.59"
$ws.Range("E14").Value = "  -6.02%  "

$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-SafeTextValue "D15" "632.43"
$ws.Range("E15").Value = "  -2.14%  "

$ws.Range("D16").Value = "3.796.73"
$ws.Range("E16").Value = "  -5.93%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-SafeTextValue "D17" "18.01"
$ws.Range("E17").Value = "  -1.22%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "65.702.09"
$ws.Range("E18").Value = "  -5.51%  "

$ws.Range("E19").Value = "  -3.49%  "

$ws.Range("D20").Value = "3.269.20"
$ws.Range("E20").Value = "  -6.37%  "

Set-SafeTextValue "D21" "11.31"
$ws.Range("E21").Value = "  -8.16%  "

Set-SafeTextValue "D22" "0.901"
$ws.Range("E22").Value = "  -4.97%  "

Set-SafeTextValue "D23" "18.29"
$ws.Range("E23").Value = "  +1.00%  "

Set-SafeTextValue "D24" "106.43"
$ws.Range("E24").Value = "  +7.61%  "

Set-SafeTextValue "D25" "4.89"
$ws.Range("E25").Value = "  -6.85%  "

$ws.Range("E26").Value = "  -7.22%  "

Set-SafeTextValue "D27" "2.67"
$ws.Range("E27").Value = "  -7.32%  "

Set-SafeTextValue "D28" "9.49"
$ws.Range("E28").Value = "  -5.49%  "

Set-SafeTextValue "D29" "8.66"
$ws.Range("E29").Value = "  -7.37%  "

Set-SafeTextValue "D30" "30.20"
$ws.Range("E30").Value = "  -7.40%  "

Set-SafeTextValue "D31" "3.94"
$ws.Range("E31").Value = "  -6.55%  "

Set-SafeTextValue "D32" "6.23"
$ws.Range("E32").Value = "  -7.46%  "

Set-SafeTextValue "D33" "11.00"
$ws.Range("E33").Value = "  -5.35%  "

$ws.Range("E34").Value = "  -4.27%  "

Set-SafeTextValue "D35" "57.57"
$ws.Range("E35").Value = "  -5.82%  "

$ws.Range("D36").Value = "3.722.94"
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("E37").Value = "  -0.05%  "

Set-SafeTextValue "D38" "523.20"
$ws.Range("E38").Value = "  -1.16%  "

Set-SafeTextValue "D39" "3.40"
$ws.Range("E39").Value = "  -4.77%  "

$ws.Range("D40").Value = "0.0₃0732"
$ws.Range("E40").Value = "  -7.07%  "

$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("E42").Value = "  -7.62%  "

$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-SafeTextValue "D43" "3.38"
$ws.Range("E43").Value = "  -3.91%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-SafeTextValue "D44" "32.82"
$ws.Range("E44").Value = "  -4.46%  "

$ws.Range("E45").Value = "  -9.88%  "

Set-SafeTextValue "D46" "3.26"
$ws.Range("E46").Value = "  -3.25%  "

$ws.Range("E47").Value = "  -6.64%  "

$ws.Range("E48").Value = "  -3.96%  "

$ws.Range("E49").Value = "  -8.57%  "

Set-SafeTextValue "D50" "0.999"
$ws.Range("E50").Value = "  -0.13%  "

Set-SafeTextValue "D51" "1.24"
$ws.Range("E51").Value = "  +1.03%  "
